# Update the answer table cells to the new values.
# The 5x5 data rows sit at table rows 1, 5, 9, 13, 17 (3 blank rows between each).
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "51÷3=17, 0"
$t.Cell(1, 2).Range.Text = "26÷4=6, 2"
$t.Cell(1, 3).Range.Text = "10÷5=2, 0"
$t.Cell(1, 4).Range.Text = "75÷2=37, 1"
$t.Cell(1, 5).Range.Text = "91÷2=45, 1"
$t.Cell(5, 1).Range.Text = "77÷4=19, 1"
$t.Cell(5, 2).Range.Text = "97÷8=12, 1"
$t.Cell(5, 3).Range.Text = "50÷7=7, 1"
$t.Cell(5, 4).Range.Text = "19÷6=3, 1"
$t.Cell(5, 5).Range.Text = "90÷3=30, 0"
$t.Cell(9, 1).Range.Text = "51÷4=12, 3"
$t.Cell(9, 2).Range.Text = "81÷6=13, 3"
$t.Cell(9, 3).Range.Text = "60÷5=12, 0"
$t.Cell(9, 4).Range.Text = "30÷2=15, 0"
$t.Cell(9, 5).Range.Text = "94÷3=31, 1"
$t.Cell(13, 1).Range.Text = "34÷5=6, 4"
$t.Cell(13, 2).Range.Text = "16÷7=2, 2"
$t.Cell(13, 3).Range.Text = "28÷5=5, 3"
$t.Cell(13, 4).Range.Text = "50÷9=5, 5"
$t.Cell(13, 5).Range.Text = "29÷5=5, 4"
$t.Cell(17, 1).Range.Text = "13÷9=1, 4"
$t.Cell(17, 2).Range.Text = "90÷3=30, 0"
$t.Cell(17, 3).Range.Text = "28÷6=4, 4"
$t.Cell(17, 4).Range.Text = "11÷7=1, 4"
$t.Cell(17, 5).Range.Text = "48÷9=5, 3"
